$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 203.5625
$ws.Range("I12").Value = 208.23077
$ws.Range("J12").Value = 183.33333
$ws.Range("K12").Value = 208.23077
$ws.Range("L12").Value = 183.33333
$ws.Range("M12").Value = -38.23077000000001
$ws.Range("N12").Value = -523.3333299999999
$ws.Range("H18").Value = 219.25
$ws.Range("I18").Value = 219.25
$ws.Range("K18").Value = 219.25
$ws.Range("M18").Value = 64.75
$ws.Range("H20").Value = 1000
$ws.Range("I20").Value = 1000
$ws.Range("K20").Value = 1000
$ws.Range("M20").Value = -770
$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("K35").Value = 1000
$ws.Range("M35").Value = -621
$ws.Range("H111").Value = 806.5
$ws.Range("I111").Value = 782.5
$ws.Range("J111").Value = 830.5
$ws.Range("K111").Value = 2347.5
$ws.Range("L111").Value = 2491.5
$ws.Range("M111").Value = 719.5
$ws.Range("N111").Value = -8625.5
$ws.Range("H127").Value = 1135.4166
$ws.Range("I127").Value = 787.5
$ws.Range("J127").Value = 1205
$ws.Range("K127").Value = 2362.5
$ws.Range("L127").Value = 3615
$ws.Range("M127").Value = 2597.5
$ws.Range("N127").Value = -13535
$ws.Range("H137").Value = 1336.7428
$ws.Range("I137").Value = 1136.625
$ws.Range("J137").Value = 1505.2632
$ws.Range("K137").Value = 3409.875
$ws.Range("L137").Value = 4515.7896
$ws.Range("M137").Value = -859.875
$ws.Range("N137").Value = -9615.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 657.7619
$ws.Range("I110").Value = 691.2778
$ws.Range("J110").Value = 456.66666
$ws.Range("K110").Value = 691.2778
$ws.Range("L110").Value = 456.66666
$ws.Range("M110").Value = 1353.7222
$ws.Range("N110").Value = -4546.66666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 26629
$ws.Range("I7").Value = 21775.5
$ws.Range("J7").Value = 36336
$ws.Range("K7").Value = 21775.5
$ws.Range("L7").Value = 36336
$ws.Range("M7").Value = -21662.5
$ws.Range("N7").Value = -36562
$ws.Range("H105").Value = 2152.75
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 1870.3334
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 1870.3334
$ws.Range("M105").Value = -1253
$ws.Range("N105").Value = -5364.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 29363.625
$ws.Range("I2").Value = 474.75
$ws.Range("J2").Value = 58252.5
$ws.Range("K2").Value = 474.75
$ws.Range("L2").Value = 58252.5
$ws.Range("M2").Value = -361.75
$ws.Range("N2").Value = -58478.5
$ws.Range("H16").Value = 1350
$ws.Range("J16").Value = 1350
$ws.Range("L16").Value = 1350
$ws.Range("N16").Value = -1924
$ws.Range("H62").Value = 500000000
$ws.Range("I62").Value = 500000000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 500000000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -499999376
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 500000000
$ws.Range("I65").Value = 500000000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 2500000000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -2499996880
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 1350
$ws.Range("J113").Value = 1350
$ws.Range("L113").Value = 1350
$ws.Range("N113").Value = -5690

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 666667400
$ws.Range("J9").Value = 666667400
$ws.Range("L9").Value = 2000002200
$ws.Range("N9").Value = -2000002648
$ws.Range("H131").Value = 940
$ws.Range("I131").Value = 441.5
$ws.Range("J131").Value = 1139.4
$ws.Range("K131").Value = 1324.5
$ws.Range("L131").Value = 3418.2
$ws.Range("M131").Value = 3715.5
$ws.Range("N131").Value = -13498.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 30000
$ws.Range("J5").Value = 30000
$ws.Range("L5").Value = 30000
$ws.Range("N5").Value = -30224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1305.5
$ws.Range("J22").Value = 1517.3334
$ws.Range("L22").Value = 1517.3334
$ws.Range("N22").Value = -2107.3334
$ws.Range("H27").Value = 1305.5
$ws.Range("J27").Value = 1517.3334
$ws.Range("L27").Value = 1517.3334
$ws.Range("N27").Value = -1731.3334
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H46").Value = 986.36365
$ws.Range("I46").Value = 992.8570999999999
$ws.Range("J46").Value = 975
$ws.Range("K46").Value = 992.8570999999999
$ws.Range("L46").Value = 975
$ws.Range("M46").Value = -804.8570999999999
$ws.Range("N46").Value = -1351
$ws.Range("H68").Value = 2032.6316
$ws.Range("I68").Value = 1707.1428
$ws.Range("J68").Value = 2944
$ws.Range("K68").Value = 1707.1428
$ws.Range("L68").Value = 2944
$ws.Range("M68").Value = -958.1428000000001
$ws.Range("N68").Value = -4442
$ws.Range("H71").Value = 2032.6316
$ws.Range("I71").Value = 1707.1428
$ws.Range("J71").Value = 2944
$ws.Range("K71").Value = 8535.714
$ws.Range("L71").Value = 14720
$ws.Range("M71").Value = -4791.714
$ws.Range("N71").Value = -22208
$ws.Range("H93").Value = 1649
$ws.Range("I93").Value = 1203.5
$ws.Range("J93").Value = 1871.75
$ws.Range("K93").Value = 1203.5
$ws.Range("L93").Value = 1871.75
$ws.Range("M93").Value = 44.5
$ws.Range("N93").Value = -4367.75
$ws.Range("H124").Value = 45788
$ws.Range("J124").Value = 45788
$ws.Range("L124").Value = 45788
$ws.Range("N124").Value = -55608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5012.5
$ws.Range("I62").Value = 4900
$ws.Range("J62").Value = 5260
$ws.Range("K62").Value = 4900
$ws.Range("L62").Value = 5260
$ws.Range("M62").Value = -4276
$ws.Range("N62").Value = -6508
$ws.Range("H65").Value = 5012.5
$ws.Range("I65").Value = 4900
$ws.Range("J65").Value = 5260
$ws.Range("K65").Value = 24500
$ws.Range("L65").Value = 26300
$ws.Range("M65").Value = -21380
$ws.Range("N65").Value = -32540
$ws.Range("H82").Value = 65291.668
$ws.Range("I82").Value = 15273
$ws.Range("J82").Value = 90301
$ws.Range("K82").Value = 15273
$ws.Range("L82").Value = 90301
$ws.Range("M82").Value = -14890
$ws.Range("N82").Value = -91067
$ws.Range("H85").Value = 65291.668
$ws.Range("I85").Value = 15273
$ws.Range("J85").Value = 90301
$ws.Range("K85").Value = 15273
$ws.Range("L85").Value = 90301
$ws.Range("M85").Value = -13947
$ws.Range("N85").Value = -92953
$ws.Range("H107").Value = 878.2963
$ws.Range("I107").Value = 718.6818
$ws.Range("J107").Value = 1580.6
$ws.Range("K107").Value = 2156.0454
$ws.Range("L107").Value = 4741.799999999999
$ws.Range("M107").Value = -236.0454
$ws.Range("N107").Value = -8581.799999999999
$ws.Range("H113").Value = 300.8889
$ws.Range("I113").Value = 259.7
$ws.Range("J113").Value = 418.57144
$ws.Range("K113").Value = 779.0999999999999
$ws.Range("L113").Value = 1255.71432
$ws.Range("M113").Value = 1390.9
$ws.Range("N113").Value = -5595.71432
$ws.Range("H132").Value = 5788.115
$ws.Range("I132").Value = 1746.4
$ws.Range("J132").Value = 19260.5
$ws.Range("K132").Value = 5239.200000000001
$ws.Range("L132").Value = 57781.5
$ws.Range("M132").Value = -2709.200000000001
$ws.Range("N132").Value = -62841.5
$ws.Range("H136").Value = 2032.1305
$ws.Range("I136").Value = 2517.2307
$ws.Range("J136").Value = 1401.5
$ws.Range("K136").Value = 7551.6921
$ws.Range("L136").Value = 4204.5
$ws.Range("M136").Value = -9304.5
